$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(44663, 73200, 121000, 122000, 121000, 122000),
    @(44664, 61600, 123500, 123000, 123500, 123000),
    @(44665, 62800, 91000, 54500, 91000, 54500)
)

$row = 43
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]

    # Match the date-format style used by the rest of column A
    $ws.Cells.Item($row, 1).NumberFormat = $ws.Cells.Item(42, 1).NumberFormat

    $row++
}
